# data2: demand node 3 demand increased to distinguish objectives 2 & 3
#
# Updates the per-period/per-commodity demand figures (columns E/F) for
# demand-node rows 3-14, the allocation figures (columns M/N/O) for row 5,
# and the capacity figures (column G) for rows 15-20. Also updates the
# active selection left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 3-8: first destination block (E: 2->6, F: 1->10) ---
foreach ($r in 3..8) {
    $ws.Cells.Item($r, 5).Value = 6   # column E
    $ws.Cells.Item($r, 6).Value = 10  # column F
}

# Row 5 additionally has its M/N/O allocation values reduced (10 -> 2)
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = 2
$ws.Range("O5").Value = 2

# --- Rows 9-14: second destination block (E: 3->5, F: 2->8) ---
foreach ($r in 9..14) {
    $ws.Cells.Item($r, 5).Value = 5   # column E
    $ws.Cells.Item($r, 6).Value = 8   # column F
}

# --- Rows 15-20: capacity column G (10 -> 50) ---
foreach ($r in 15..20) {
    $ws.Cells.Item($r, 7).Value = 50  # column G
}

# --- Restore the author's active selection ---
$ws.Range("P6").Select()
